# Update StructureDefinition-episode-disease-stage-code.xlsx
# (IG build regenerated this spreadsheet export: version bump, new build
# date, publisher/jurisdiction metadata, and the profile's own
# short/definition text in the Elements sheet.)

$wb = $excel.ActiveWorkbook

# --- "Metadata" sheet -------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")

# Version: 5.0.0 -> 6.0.0
$meta.Range("B3").Value = "6.0.0"

# Date: regenerated build timestamp
$meta.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher row previously had an empty value and was immediately followed
# by a duplicated "Contact" / "No display for ContactDetail" row. The new
# build fills in the Publisher value and replaces the duplicate row with a
# new Jurisdiction row, so everything from row 9 down collapses by one row.
$meta.Range("B9").Value = "Alvearie Team"

$meta.Range("A10").Value = "Jurisdiction"
$meta.Range("B10").Value = "United States of America"

$meta.Range("A11").Value = "Description"
$meta.Range("B11").Value = "Disease stage code for the episode of care"

$meta.Range("A12").Value = "Purpose"
$meta.Range("B12").Value = ""

$meta.Range("A13").Value = "Copyright"
$meta.Range("B13").Value = ""

$meta.Range("A14").Value = "FHIR Version"
$meta.Range("B14").Value = "4.0.1"

$meta.Range("A15").Value = "Kind"
$meta.Range("B15").Value = "complex-type"

$meta.Range("A16").Value = "Type"
$meta.Range("B16").Value = "Extension"

$meta.Range("A17").Value = "Base Definition"
$meta.Range("B17").Value = "http://hl7.org/fhir/StructureDefinition/Extension"

$meta.Range("A18").Value = "Abstract"
# Plain Value assignment of the bare word "false" auto-converts to the
# Boolean FALSE in Excel; round-trip it through a formula + paste-as-values
# so the cell keeps storing the literal text "false" (as the source sheet
# does for every other property value here).
$meta.Range("B18").Formula = "=""false"""
$meta.Range("B18").Copy()
$meta.Range("B18").PasteSpecial(-4163)

$meta.Range("A19").Value = "Derivation"
$meta.Range("B19").Value = "constraint"

$meta.Range("A20").Value = "Context"
$meta.Range("B20").Value = "element:Element"

# The old sheet had one extra trailing row (21); delete it so the sheet
# ends at row 20, matching the new dimension A1:B20.
$meta.Rows.Item(21).Delete()

# --- "Elements" sheet ---------------------------------------------------
# The root Extension row's Short/Definition columns are updated from the
# generic "Extension" / "An Extension" placeholders to the profile-specific
# title/description text.
$elements = $wb.Worksheets.Item("Elements")
$elements.Range("K2").Value = "Episode Disease Stage Code"
$elements.Range("L2").Value = "Disease stage code for the episode of care"
